$d = $word.ActiveDocument

# The "Objects & Design," course in the coursework list is being dropped
# (it is effectively folded into the "Design & Analysis of Algorithms..."
# course that already follows it in the same sentence).
$removed = $d.Content.Find.Execute("Objects & Design, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
if (-not $removed) {
    throw "Could not find 'Objects & Design, ' to remove"
}

# "Operating System Design" is added as the first course listed after the
# "Coursework:" label.
$added = $d.Content.Find.Execute("Coursework: ", $true, $false, $false, $false, $false, $true, 1, $false, "Coursework: Operating System Design, ", 2)
if (-not $added) {
    throw "Could not find 'Coursework: ' to update"
}
